# Update "Forecast Comparison" sheet with correct forecast output:
#  - insert a new "Week_Start_Date" column after "Week" (new col B)
#  - shorten the week labels in col A (W01 -> W1 ... W09 -> W9, W10.. unchanged)
#  - re-type the is_holiday_week column (now col J) as boolean

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Forecast Comparison")

# Insert a new column before the old column B (ASIN), pushing ASIN.. onward
# one column to the right (B..I -> C..J).
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Cells.Item(1, 2).Value = "Week_Start_Date"

# Per-week data: short week label, week-start date (kept as literal text,
# not auto-converted to a date serial), and the five forecast numbers that
# used to live in columns C:G and now live in D:H.
$weeks = @(
    @{ Row = 2;  Label = "W1";  Date = "2025-01-05"; My = 5; Mean = 6;  P70 = 6;  P80 = 11; P90 = 19 },
    @{ Row = 3;  Label = "W2";  Date = "2025-01-12"; My = 6; Mean = 7;  P70 = 7;  P80 = 12; P90 = 20 },
    @{ Row = 4;  Label = "W3";  Date = "2025-01-19"; My = 7; Mean = 8;  P70 = 10; P80 = 14; P90 = 23 },
    @{ Row = 5;  Label = "W4";  Date = "2025-01-26"; My = 8; Mean = 9;  P70 = 11; P80 = 16; P90 = 24 },
    @{ Row = 6;  Label = "W5";  Date = "2025-02-02"; My = 8; Mean = 9;  P70 = 10; P80 = 15; P90 = 24 },
    @{ Row = 7;  Label = "W6";  Date = "2025-02-09"; My = 8; Mean = 9;  P70 = 10; P80 = 15; P90 = 24 },
    @{ Row = 8;  Label = "W7";  Date = "2025-02-16"; My = 8; Mean = 9;  P70 = 10; P80 = 15; P90 = 25 },
    @{ Row = 9;  Label = "W8";  Date = "2025-02-23"; My = 8; Mean = 9;  P70 = 10; P80 = 15; P90 = 25 },
    @{ Row = 10; Label = "W9";  Date = "2025-03-02"; My = 7; Mean = 8;  P70 = 9;  P80 = 14; P90 = 24 },
    @{ Row = 11; Label = "W10"; Date = "2025-03-09"; My = 8; Mean = 9;  P70 = 9;  P80 = 15; P90 = 24 },
    @{ Row = 12; Label = "W11"; Date = "2025-03-16"; My = 8; Mean = 9;  P70 = 10; P80 = 16; P90 = 27 },
    @{ Row = 13; Label = "W12"; Date = "2025-03-23"; My = 8; Mean = 9;  P70 = 10; P80 = 16; P90 = 27 },
    @{ Row = 14; Label = "W13"; Date = "2025-03-30"; My = 8; Mean = 9;  P70 = 10; P80 = 15; P90 = 25 },
    @{ Row = 15; Label = "W14"; Date = "2025-04-06"; My = 8; Mean = 9;  P70 = 9;  P80 = 15; P90 = 27 },
    @{ Row = 16; Label = "W15"; Date = "2025-04-13"; My = 8; Mean = 10; P70 = 10; P80 = 16; P90 = 28 },
    @{ Row = 17; Label = "W16"; Date = "2025-04-20"; My = 8; Mean = 10; P70 = 10; P80 = 16; P90 = 28 }
)

foreach ($w in $weeks) {
    $r = $w.Row

    # Col A: short week label (strip leading zero from W01..W09).
    $ws.Cells.Item($r, 1).Value = $w.Label

    # Col B: week start date, forced to stay plain text (Excel would
    # otherwise auto-convert an ISO-looking string into a date serial).
    $cellB = $ws.Cells.Item($r, 2)
    $cellB.NumberFormat = "@"
    $cellB.Value = $w.Date

    # Cols D:H: shifted forecast numbers (ASIN itself, col C, is untouched
    # by the insert and keeps its original "B07FW85VFT" value).
    $ws.Cells.Item($r, 4).Value = $w.My
    $ws.Cells.Item($r, 5).Value = $w.Mean
    $ws.Cells.Item($r, 6).Value = $w.P70
    $ws.Cells.Item($r, 7).Value = $w.P80
    $ws.Cells.Item($r, 8).Value = $w.P90

    # Col J (was I, shifted by the insert): is_holiday_week, now typed as
    # a genuine boolean cell rather than a 0/1 number.
    $ws.Cells.Item($r, 10).Value = $false
}
